# Actualización automática 2025-12-01 16:30:09
# Updates sales figures for HIDALGO HIDALGO PEDRO GUSTAVO row (row 22)
# on the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, plus their
# corresponding "X de 21" completion counters / totals on row 23.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 22 - HIDALGO HIDALGO PEDRO GUSTAVO's sales by product group
$wsGrupo.Range("D22").Value = 457.92
$wsGrupo.Range("H22").Value = 2847.34
$wsGrupo.Range("M22").Value = 5082.56

# Row 23 - completion counters ("X de 21") for the columns that now have
# a non-zero entry
$wsGrupo.Range("D23").Value = "1 de 21"
$wsGrupo.Range("H23").Value = "1 de 21"
$wsGrupo.Range("M23").Value = "1 de 21"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 22 - HIDALGO HIDALGO PEDRO GUSTAVO's December sales value
$wsMensual.Range("F22").Value = 8387.82

# Row 23 - column total for December
$wsMensual.Range("F23").Value = 8387.82
